$p = $ppt.ActivePresentation

# --- Slide 12 ("Ausblick"): remove the leftover "Untertitel 2" caption
#     textbox ("Abbildung: Typen von Zusammenfassungen nach Gambhir/ Gupta.")
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item("Untertitel 2").Delete()

# --- Slide 13 ("Arbeit anmelden"): mark the last bullet with a red "(!)"
$s13 = $p.Slides.Item(13)
$contentShape = $s13.Shapes.Item(1)
$contentRange = $contentShape.TextFrame.TextRange
$lastParaIndex = $contentRange.Paragraphs().Count
$lastPara = $contentRange.Paragraphs($lastParaIndex)

$firstRun = $lastPara.Runs(1)
$baseLen = $firstRun.Text.Length
$firstRun.Text = $firstRun.Text + " (!)"

$exclamStart = $baseLen + 3
$exclam = $lastPara.Characters($exclamStart, 1)
$exclam.Font.Color.RGB = 3277010
